$d = $word.ActiveDocument

# Section headings whose body paragraphs (between this heading and the
# next "Normal"-styled separator paragraph) must be struck through.
$targetHeadings = @("PROMPT 12: Component Fixtures", "PROMPT 13: Advanced Reporting")

$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $styleName = $p.Range.ParagraphStyle.NameLocal
    $text = $p.Range.Text

    if ($styleName -eq "Heading 3") {
        $trimmed = $text.Trim()
        $isTarget = $false
        foreach ($h in $targetHeadings) {
            if ($trimmed -eq $h) {
                $isTarget = $true
            }
        }

        if ($isTarget) {
            $j = $i + 1
            while ($j -le $count) {
                $bodyPara = $d.Paragraphs($j)
                $bodyStyle = $bodyPara.Range.ParagraphStyle.NameLocal
                if ($bodyStyle -eq "Normal") {
                    break
                }
                $bodyPara.Range.Font.StrikeThrough = $true
                $j = $j + 1
            }
        }
    }
}
